$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a number by
# Excels .Value setter (losing formatting like trailing zeros / sci notation)
# are briefly forced to Text format, then restored to the default "Normal"
# style so no stray numFmt/style is left on the cell (matches source: all
# data cells in column D/E use the workbook default style).

$ws.Range('D2').Value = '65.194.30'
$ws.Range('E2').Value = '  +2.92%  '
$ws.Range('D3').Value = '3.401.83'
$ws.Range('E3').Value = '  +2.57%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '562.23'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '174.90'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.77%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.626'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.07%  '
$ws.Range('D8').Value = '3.391.76'
$ws.Range('E8').Value = '  +2.73%  '
$ws.Range('E9').Value = '  -0.07%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.168'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +12.59%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.632'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.21%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '54.75'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.14%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000278'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +5.72%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.15'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.09%  '
$ws.Range('D15').Value = '3.942.95'
$ws.Range('E15').Value = '  +2.09%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '18.39'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +4.06%  '
$ws.Range('D17').Value = '3.400.40'
$ws.Range('E17').Value = '  +2.61%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.119'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.21%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.92'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.62%  '
$ws.Range('D20').Value = '65.130.00'
$ws.Range('E20').Value = '  +2.75%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.996'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.18%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '474.18'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +17.39%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.96'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +16.36%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.15'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.72%  '
$ws.Range('B25').Value = 'InternetComputer(DFINITY)'
$ws.Range('C25').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '13.84'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +5.37%  '
$ws.Range('B26').Value = 'Litecoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '86.73'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +5.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.90'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.34%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.89'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.23%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.85'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.43%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '30.74'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.93%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.70'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.67%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.55'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.54%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '585.10'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.40%  '
$ws.Range('E34').Value = '  +3.43%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '60.42'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.35%  '
$ws.Range('E36').Value = '  +0.19%  '
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.141'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.84%  '
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.53'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.93%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '35.99'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.93%  '
$ws.Range('D40').Value = '0.0₃0750'
$ws.Range('E40').Value = '  +2.39%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.374'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.54%  '
$ws.Range('D42').Value = '3.112.75'
$ws.Range('E42').Value = '  -0.91%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.999'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.21%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.86'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.55%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.52'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.77%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0416'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.75%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.22'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.22%  '
$ws.Range('E48').Value = '  +5.31%  '
$ws.Range('E49').Value = '  -1.45%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '136.95'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.40%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.35'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.59%  '
